$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "HTTP" results row (row 7) with the new measured values.
$ws.Range("B7").Value = 566.24
$ws.Range("C7").Value = 36.14
$ws.Range("D7").Value = 55946.98
$ws.Range("E7").Value = 3105.75
$ws.Range("F7").Value = 47887.27
$ws.Range("G7").Value = 5395.2
$ws.Range("H7").Value = 1036.96
$ws.Range("I7").Value = 4.67
$ws.Range("J7").Value = 3.0099900000000002
$ws.Range("K7").Value = 1.02725
$ws.Range("L7").Value = 1.0081899999999999
$ws.Range("M7").Value = 1.0158

# Keep the active selection on the cell that was last edited, matching the
# saved workbook's view state.
$null = $ws.Range("C7").Select()
